# Adding the changes we made on may 9th
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new blank rows above the existing data block (rows 2-8),
# which pushes the current data rows (2-21) down to rows 9-28.
$ws.Range("A2:H8").EntireRow.Insert()

# The inserted rows pick up the header row's formatting by default;
# strip it so the new data rows look like ordinary (unstyled) data rows.
$ws.Range("A2:H8").ClearFormats()

# Fill in the 7 newly inserted rows with the new sensor readings.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "falling"
$ws.Range("C2").Value = 1.19674801826477
$ws.Range("D2").Value = 1.663910716772079
$ws.Range("E2").Value = 2.702408194541931
$ws.Range("F2").Value = -0.1724167168140411
$ws.Range("G2").Value = -0.3089450895786285
$ws.Range("H2").Value = 0.9990701079368592

$ws.Range("A3").Value = 100
$ws.Range("B3").Value = "falling"
$ws.Range("C3").Value = -0.4646213054656981
$ws.Range("D3").Value = 2.042550325393677
$ws.Range("E3").Value = 2.406269252300262
$ws.Range("F3").Value = -0.1458440721035003
$ws.Range("G3").Value = -0.0731511116027832
$ws.Range("H3").Value = -0.1902845203876495

$ws.Range("A4").Value = 200
$ws.Range("B4").Value = "falling"
$ws.Range("C4").Value = -0.3934619426727292
$ws.Range("D4").Value = 1.991465017199517
$ws.Range("E4").Value = 1.958218067884445
$ws.Range("F4").Value = 0.0277943685650825
$ws.Range("G4").Value = -0.0499382354319095
$ws.Range("H4").Value = 0.04505131021142

$ws.Range("A5").Value = 300
$ws.Range("B5").Value = "falling"
$ws.Range("C5").Value = -0.06526184082031269
$ws.Range("D5").Value = 1.843156695365906
$ws.Range("E5").Value = 2.04642915725708
$ws.Range("F5").Value = -0.042302418500185
$ws.Range("G5").Value = -0.052381694316864
$ws.Range("H5").Value = -0.0262672062963247

$ws.Range("A6").Value = 400
$ws.Range("B6").Value = "falling"
$ws.Range("C6").Value = -0.2364732027053833
$ws.Range("D6").Value = 1.819270551204681
$ws.Range("E6").Value = 2.093152940273285
$ws.Range("F6").Value = 0.07803803682327271
$ws.Range("G6").Value = -0.0073303831741213
$ws.Range("H6").Value = 0.0215329993516206

$ws.Range("A7").Value = 500
$ws.Range("B7").Value = "falling"
$ws.Range("C7").Value = -0.3396859169006348
$ws.Range("D7").Value = 1.832332909107209
$ws.Range("E7").Value = 2.310090780258179
$ws.Range("F7").Value = -0.0259617734700441
$ws.Range("G7").Value = -0.0493273697793483
$ws.Range("H7").Value = -0.0320704244077205

$ws.Range("A8").Value = 600
$ws.Range("B8").Value = "falling"
$ws.Range("C8").Value = -0.2427999973297116
$ws.Range("D8").Value = 1.836586102843285
$ws.Range("E8").Value = 2.258781224489212
$ws.Range("F8").Value = -0.117286130785942
$ws.Range("G8").Value = -0.0560468845069408
$ws.Range("H8").Value = -0.0174096599221229

# Append 3 brand-new rows after the shifted data (now ending at row 28).
$ws.Range("A29").Value = 2700
$ws.Range("B29").Value = "falling"
$ws.Range("C29").Value = -0.2655735015869125
$ws.Range("D29").Value = 2.233672142028808
$ws.Range("E29").Value = 0.9439086914062514
$ws.Range("F29").Value = 0.1461495161056518
$ws.Range("G29").Value = 0.1336267739534378
$ws.Range("H29").Value = -0.1892155110836029

$ws.Range("A30").Value = 2800
$ws.Range("B30").Value = "falling"
$ws.Range("C30").Value = 0.07992589473724532
$ws.Range("D30").Value = 1.95888604223728
$ws.Range("E30").Value = 1.25704461336136
$ws.Range("F30").Value = 0.117286130785942
$ws.Range("G30").Value = 0.7583891749382019
$ws.Range("H30").Value = 0.07605272531509399

$ws.Range("A31").Value = 2900
$ws.Range("B31").Value = "falling"
$ws.Range("C31").Value = 0.3567421436309829
$ws.Range("D31").Value = 2.357963830232623
$ws.Range("E31").Value = 1.160924613475799
$ws.Range("F31").Value = 0.1061378344893455
$ws.Range("G31").Value = 0.2086104750633239
$ws.Range("H31").Value = -0.1314887404441833

Write-Output "applied 10 new rows (A2:H8 inserted, A29:H31 appended)"
